$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New section: "Array Sum" parallel timing tables (rows 16-20 and 22-27)
# ---------------------------------------------------------------------------

# --- First table (row 16 header, rows 17-20 data) ---------------------------
$ws.Range("B16").Value = "n"
$ws.Range("C16").Value = "t(sec)"
$ws.Range("D16").Value = "1 thread"
$ws.Range("E16").Value = "2 thread"
$ws.Range("F16").Value = "4 thread"
$ws.Range("G16").Value = "8 thread"
$ws.Range("D16:G16").NumberFormat = "0.00000"

$ws.Range("B17").Value = 10000
$ws.Range("C17").Value = 0.027
$ws.Range("D17").Value = 0.003
$ws.Range("E17").Value = 0.004
$ws.Range("F17").Value = 0.004
$ws.Range("G17").Value = 0.004

$ws.Range("B18").Value = 100000
$ws.Range("C18").Value = 0.237
$ws.Range("D18").Value = 0.016
$ws.Range("E18").Value = 0.015
$ws.Range("F18").Value = 0.014
$ws.Range("G18").Value = 0.015

$ws.Range("B19").Value = 1000000
$ws.Range("C19").Value = 2.497
$ws.Range("D19").Value = 0.124
$ws.Range("E19").Value = 0.124
$ws.Range("F19").Value = 0.124
$ws.Range("G19").Value = 0.124

$ws.Range("B20").Value = 10000000
$ws.Range("C20").Value = 24.793
$ws.Range("D20").Value = 1.217
$ws.Range("E20").Value = 1.209
$ws.Range("F20").Value = 1.206
$ws.Range("G20").Value = 1.204

# --- Second table (row 22 title, row 23 header, rows 24-27 data) -----------
$ws.Range("C22").Value = "t((sec)"
$ws.Range("B23").Value = "n"
$ws.Range("C23").Value = "1 thread"
$ws.Range("D23").Value = "2 thread"
$ws.Range("E23").Value = "4 thread"
$ws.Range("F23").Value = "8 thread"

$ws.Range("B24").Value = 10000
$ws.Range("C24").Value = 0.003
$ws.Range("D24").Value = 0.004
$ws.Range("E24").Value = 0.004
$ws.Range("F24").Value = 0.004

$ws.Range("B25").Value = 100000
$ws.Range("C25").Value = 0.016
$ws.Range("D25").Value = 0.015
$ws.Range("E25").Value = 0.014
$ws.Range("F25").Value = 0.015

$ws.Range("B26").Value = 1000000
$ws.Range("C26").Value = 0.124
$ws.Range("D26").Value = 0.124
$ws.Range("E26").Value = 0.124
$ws.Range("F26").Value = 0.124

$ws.Range("B27").Value = 10000000
$ws.Range("C27").Value = 1.217
$ws.Range("D27").Value = 1.209
$ws.Range("E27").Value = 1.206
$ws.Range("F27").Value = 1.204

# Give the second table its thin-box border + number format (header row
# keeps the "0.00000" numeric format applied, like the first table's header).
$ws.Range("B22:F27").Borders.LineStyle = 1
$ws.Range("C23:F23").NumberFormat = "0.00000"

# Selection left where the author last clicked.
$ws.Range("J18").Select() | Out-Null
